# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets. Both sheets contain identical data
# tables, so the same F-column updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2" = 1325
    "F3" = 1804
    "F4" = 124
    "F6" = 6286
    "F7" = 142
    "F8" = 109
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
